$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "ICs" sheet: insert a new row for a 7812 (TO-220) linear regulator,
#    right after the existing 7805 regulator row (old row 25 -> still 25).
# ---------------------------------------------------------------------------
$ics = $wb.Worksheets.Item("ICs")

$ics.Rows.Item(26).Insert()

$ics.Cells.Item(26,1).Value  = "7812"
$ics.Cells.Item(26,2).Value  = "TO-220"
$ics.Cells.Item(26,3).Value  = "TO-220"
$ics.Cells.Item(26,4).Value  = "L7812CV"
$ics.Cells.Item(26,5).Value  = "STMicroelectronics"
$ics.Cells.Item(26,6).Value  = "L7812CV"
$ics.Cells.Item(26,7).Value  = "Digi-Key"
$ics.Cells.Item(26,8).Value  = "497-1452-5-ND"
$ics.Cells.Item(26,9).Value  = "Fitted"
$ics.Cells.Item(26,10).Value = "Generic"
$ics.Cells.Item(26,13).Value = "12V 1.5A linear regulator, 14-35V input"
$ics.Cells.Item(26,15).Value = "PMIC"
$ics.Cells.Item(26,16).Value = "y"

# ---------------------------------------------------------------------------
# 2) "Modules" sheet: insert a new row for the plain Feather RP2040 board
#    (before the Feather RP2040 RFM95 row).
# ---------------------------------------------------------------------------
$modules = $wb.Worksheets.Item("Modules")

$modules.Rows.Item(6).Insert()

$modules.Cells.Item(6,1).Value  = "Feather RP2040"
$modules.Cells.Item(6,2).Value  = "Feather RP2040"
$modules.Cells.Item(6,3).Value  = "Feather RP2040"
$modules.Cells.Item(6,4).Value  = "Feather RP2040"
$modules.Cells.Item(6,5).Value  = "Adafruit"
$modules.Cells.Item(6,6).Value  = "4884"
$modules.Cells.Item(6,7).Value  = "Digi-Key"
$modules.Cells.Item(6,8).Value  = "1528-4884-ND"
$modules.Cells.Item(6,9).Value  = "Fitted"
$modules.Cells.Item(6,10).Value = "Non Generic"
$modules.Cells.Item(6,13).Value = "Feather RP2040"
$modules.Cells.Item(6,15).Value = "y"

# ---------------------------------------------------------------------------
# 3) "ICs" sheet: append 5 new rows for Xilinx Artix-7 FPGAs (XC7AxxT-1FTG256
#    family) at the bottom of the sheet.
# ---------------------------------------------------------------------------
$fpgaRows = @(
    @{ D = "XC7A15T-1TFG256C"; F = "XC7A15T-1FTG256C"; H = "122-1930-ND" },
    @{ D = "XC7A35T-1TFG256C"; F = "XC7A35T-1FTG256C"; H = "122-1910-ND" },
    @{ D = "XC7A50T-1TFG256C"; F = "XC7A50T-1FTG256C"; H = "122-1916-ND" },
    @{ D = "XC7A75T-1TFG256C"; F = "XC7A75T-1FTG256C"; H = "XC7A75T-1FTG256C-ND" },
    @{ D = "XC7A100T-1TFG256C"; F = "XC7A100T-1FTG256C"; H = "122-1882-ND" }
)

$startRow = 47
for ($i = 0; $i -lt $fpgaRows.Count; $i++) {
    $r = $startRow + $i
    $row = $fpgaRows[$i]

    $ics.Cells.Item($r,1).Value  = "XC7AxxT-1FTG256"
    $ics.Cells.Item($r,2).Value  = "XC7AxxT-1FTG256"
    $ics.Cells.Item($r,3).Value  = "FTG256"
    $ics.Cells.Item($r,4).Value  = $row.D
    $ics.Cells.Item($r,5).Value  = "Xilinx"
    $ics.Cells.Item($r,6).Value  = $row.F
    $ics.Cells.Item($r,7).Value  = "Digi-Key"
    $ics.Cells.Item($r,8).Value  = $row.H
    $ics.Cells.Item($r,9).Value  = "Fitted"
    $ics.Cells.Item($r,10).Value = "Non Generic"
    $ics.Cells.Item($r,13).Value = "Artix-7 FPGA"
    $ics.Cells.Item($r,15).Value = "FPGA"
    $ics.Cells.Item($r,16).Value = "y"
}

# The new Modules row pushes the longest "Description" text further down the
# column, and Excel re-derives the best-fit width for column M (Description)
# accordingly once data entry is complete.
$modules.Columns.Item(13).ColumnWidth = 41.3

# ---------------------------------------------------------------------------
# 4) Final view state: ICs becomes the active sheet/tab, with a specific
#    selected cell; Modules keeps a lingering selection further down too.
# ---------------------------------------------------------------------------
$modules.Range("N20").Select()

$ics.Activate()
$ics.Range("J31").Select()
